$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.826.81'
$ws.Range("E2").Value = '  +1.94%  '
$ws.Range("D3").Value = '3.478.16'
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.85'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.10'
$ws.Range("E6").Value = '  +3.92%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.65'
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.399'
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("D12").Value = '4.076.82'
$ws.Range("E12").Value = '  +2.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.90'
$ws.Range("E13").Value = '  +5.11%  '
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '3.483.35'
$ws.Range("E15").Value = '  +2.86%  '
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '62.946.28'
$ws.Range("E17").Value = '  +2.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.33'
$ws.Range("E18").Value = '  +2.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.42'
$ws.Range("E19").Value = '  +5.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.34'
$ws.Range("E20").Value = '  +3.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.71'
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.565'
$ws.Range("E22").Value = '  +1.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.15'
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = '3.625.38'
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000116'
$ws.Range("E26").Value = '  +2.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.179'
$ws.Range("E27").Value = '  -7.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.66'
$ws.Range("E28").Value = '  +4.86%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.23'
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.13'
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.40'
$ws.Range("E32").Value = '  +1.85%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.75'
$ws.Range("E34").Value = '  +1.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.11'
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.27'
$ws.Range("E36").Value = '  +3.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '31.74'
$ws.Range("E37").Value = '  +22.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '171.81'
$ws.Range("E38").Value = '  +2.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.57'
$ws.Range("E39").Value = '  +6.48%  '
$ws.Range("D40").Value = '3.521.94'
$ws.Range("E40").Value = '  +2.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0769'
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.807'
$ws.Range("E42").Value = '  +3.23%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.22'
$ws.Range("E43").Value = '  -0.76%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.46'
$ws.Range("E44").Value = '  +0.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.71'
$ws.Range("E45").Value = '  +2.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.19'
$ws.Range("E46").Value = '  +3.05%  '
$ws.Range("D47").Value = '2.602.23'
$ws.Range("E47").Value = '  +5.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.53'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.25'
$ws.Range("E49").Value = '  +8.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.79'
$ws.Range("E50").Value = '  +1.47%  '
$ws.Range("E51").Value = '  +0.14%  '
